$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value2 = 676.1  # H2: 625.5454999999999 -> 676.1
$ws.Cells.Item(2,10).Value2 = 1079.4  # J2: 919.5 -> 1079.4
$ws.Cells.Item(2,12).Value2 = 1079.4  # L2: 919.5 -> 1079.4
$ws.Cells.Item(2,14).Value2 = -1305.4  # N2: -1145.5 -> -1305.4
$ws.Cells.Item(19,8).Value2 = 1124.2142  # H19: 1126.9286 -> 1124.2142
$ws.Cells.Item(19,9).Value2 = 220  # I19: 229.5 -> 220
$ws.Cells.Item(19,11).Value2 = 220  # K19: 229.5 -> 220
$ws.Cells.Item(19,13).Value2 = -45  # M19: -54.5 -> -45
$ws.Cells.Item(38,8).Value2 = 1728.5454  # H38: 3040.75 -> 1728.5454
$ws.Cells.Item(38,9).Value2 = 64.875  # I38: 62.25 -> 64.875
$ws.Cells.Item(38,10).Value2 = 6165  # J38: 8997.75 -> 6165
$ws.Cells.Item(38,11).Value2 = 194.625  # K38: 186.75 -> 194.625
$ws.Cells.Item(38,12).Value2 = 18495  # L38: 26993.25 -> 18495
$ws.Cells.Item(38,13).Value2 = 177.375  # M38: 185.25 -> 177.375
$ws.Cells.Item(38,14).Value2 = -19239  # N38: -27737.25 -> -19239
$ws.Cells.Item(98,8).Value2 = 4496.636  # H98: 4547.6 -> 4496.636
$ws.Cells.Item(98,9).Value2 = 4496.636  # I98: 4547.6 -> 4496.636
$ws.Cells.Item(98,11).Value2 = 4496.636  # K98: 4547.6 -> 4496.636
$ws.Cells.Item(98,13).Value2 = -2998.636  # M98: -3049.6 -> -2998.636
$ws.Cells.Item(100,8).Value2 = 914.7143  # H100: 917.2143 -> 914.7143
$ws.Cells.Item(100,9).Value2 = 915.75  # I100: 920.125 -> 915.75
$ws.Cells.Item(100,11).Value2 = 915.75  # K100: 920.125 -> 915.75
$ws.Cells.Item(100,13).Value2 = -374.75  # M100: -379.125 -> -374.75
$ws.Cells.Item(122,8).Value2 = 4496.636  # H122: 4547.6 -> 4496.636
$ws.Cells.Item(122,9).Value2 = 4496.636  # I122: 4547.6 -> 4496.636
$ws.Cells.Item(122,11).Value2 = 13489.908  # K122: 13642.8 -> 13489.908
$ws.Cells.Item(122,13).Value2 = -11039.908  # M122: -11192.8 -> -11039.908
$ws.Cells.Item(132,8).Value2 = 38145.41  # H132: 39542.15 -> 38145.41
$ws.Cells.Item(132,9).Value2 = 40999.29  # I132: 41001.848 -> 40999.29
$ws.Cells.Item(132,10).Value2 = 1045  # J132: 1590 -> 1045
$ws.Cells.Item(132,11).Value2 = 122997.87  # K132: 123005.544 -> 122997.87
$ws.Cells.Item(132,12).Value2 = 3135  # L132: 4770 -> 3135
$ws.Cells.Item(132,13).Value2 = -120467.87  # M132: -120475.544 -> -120467.87
$ws.Cells.Item(132,14).Value2 = -8195  # N132: -9830 -> -8195
$ws.Cells.Item(134,8).Value2 = 199898.5  # H134: 199948.6 -> 199898.5
$ws.Cells.Item(134,10).Value2 = 199898.5  # J134: 199948.6 -> 199898.5
$ws.Cells.Item(134,12).Value2 = 199898.5  # L134: 199948.6 -> 199898.5
$ws.Cells.Item(134,14).Value2 = -210038.5  # N134: -210088.6 -> -210038.5
$ws.Cells.Item(138,8).Value2 = 2666.9092  # H138: 2624.6155 -> 2666.9092
$ws.Cells.Item(138,10).Value2 = 3979.4644  # J138: 4049 -> 3979.4644
$ws.Cells.Item(138,12).Value2 = 11938.3932  # L138: 12147 -> 11938.3932
$ws.Cells.Item(138,14).Value2 = -22218.3932  # N138: -22427 -> -22218.3932

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2,8).Value2 = 3328245  # H2: 3536229.8 -> 3328245
$ws.Cells.Item(2,9).Value2 = 5656372.5  # I2: 6284804 -> 5656372.5
$ws.Cells.Item(2,11).Value2 = 5656372.5  # K2: 6284804 -> 5656372.5
$ws.Cells.Item(2,13).Value2 = -5656259.5  # M2: -6284691 -> -5656259.5
$ws.Cells.Item(24,8).Value2 = 28355  # H24: 0 -> 28355
$ws.Cells.Item(24,10).Value2 = 28355  # J24: 0 -> 28355
$ws.Cells.Item(24,12).Value2 = 28355  # L24: 0 -> 28355
$ws.Cells.Item(24,14).Value2 = -29103  # N24: None -> -29103
$ws.Cells.Item(32,8).Value2 = 6427.4688  # H32: 6494.8105 -> 6427.4688
$ws.Cells.Item(32,9).Value2 = 3794.9746  # I32: 3843.2437 -> 3794.9746
$ws.Cells.Item(32,11).Value2 = 3794.9746  # K32: 3843.2437 -> 3794.9746
$ws.Cells.Item(32,13).Value2 = -3507.9746  # M32: -3556.2437 -> -3507.9746
$ws.Cells.Item(61,8).Value2 = 17110.334  # H61: 21199.143 -> 17110.334
$ws.Cells.Item(61,9).Value2 = 18624.75  # I61: 23899.834 -> 18624.75
$ws.Cells.Item(61,11).Value2 = 18624.75  # K61: 23899.834 -> 18624.75
$ws.Cells.Item(61,13).Value2 = -18412.75  # M61: -23687.834 -> -18412.75
$ws.Cells.Item(96,8).Value2 = 32419.25  # H96: 32835.75 -> 32419.25
$ws.Cells.Item(96,10).Value2 = 32419.25  # J96: 32835.75 -> 32419.25
$ws.Cells.Item(96,12).Value2 = 32419.25  # L96: 32835.75 -> 32419.25
$ws.Cells.Item(96,14).Value2 = -37911.25  # N96: -38327.75 -> -37911.25
$ws.Cells.Item(100,8).Value2 = 28355  # H100: 0 -> 28355
$ws.Cells.Item(100,10).Value2 = 28355  # J100: 0 -> 28355
$ws.Cells.Item(100,12).Value2 = 28355  # L100: 0 -> 28355
$ws.Cells.Item(100,14).Value2 = -30519  # N100: None -> -30519
$ws.Cells.Item(110,8).Value2 = 2527049.2  # H110: 3088198 -> 2527049.2
$ws.Cells.Item(110,9).Value2 = 3473919  # I110: 3969940.2 -> 3473919
$ws.Cells.Item(110,10).Value2 = 2063  # J110: 2099.5 -> 2063
$ws.Cells.Item(110,11).Value2 = 3473919  # K110: 3969940.2 -> 3473919
$ws.Cells.Item(110,12).Value2 = 2063  # L110: 2099.5 -> 2063
$ws.Cells.Item(110,13).Value2 = -3471874  # M110: -3967895.2 -> -3471874
$ws.Cells.Item(110,14).Value2 = -6153  # N110: -6189.5 -> -6153
$ws.Cells.Item(116,8).Value2 = 3328245  # H116: 3536229.8 -> 3328245
$ws.Cells.Item(116,9).Value2 = 5656372.5  # I116: 6284804 -> 5656372.5
$ws.Cells.Item(116,11).Value2 = 5656372.5  # K116: 6284804 -> 5656372.5
$ws.Cells.Item(116,13).Value2 = -5654078.5  # M116: -6282510 -> -5654078.5
$ws.Cells.Item(122,8).Value2 = 2229584.8  # H122: 2373367.8 -> 2229584.8
$ws.Cells.Item(122,9).Value2 = 2107399.8  # I122: 2290569.8 -> 2107399.8
$ws.Cells.Item(122,11).Value2 = 6322199.399999999  # K122: 6871709.399999999 -> 6322199.399999999
$ws.Cells.Item(122,13).Value2 = -6319749.399999999  # M122: -6869259.399999999 -> -6319749.399999999
$ws.Cells.Item(132,8).Value2 = 5712.4  # H132: 4326.5713 -> 5712.4
$ws.Cells.Item(132,9).Value2 = 7033.1577  # I132: 4594.5483 -> 7033.1577
$ws.Cells.Item(132,10).Value2 = 4144  # J132: 3865.0557 -> 4144
$ws.Cells.Item(132,11).Value2 = 21099.4731  # K132: 13783.6449 -> 21099.4731
$ws.Cells.Item(132,12).Value2 = 12432  # L132: 11595.1671 -> 12432
$ws.Cells.Item(132,13).Value2 = -18569.4731  # M132: -11253.6449 -> -18569.4731
$ws.Cells.Item(132,14).Value2 = -17492  # N132: -16655.1671 -> -17492
$ws.Cells.Item(136,8).Value2 = 17110.334  # H136: 21199.143 -> 17110.334
$ws.Cells.Item(136,9).Value2 = 18624.75  # I136: 23899.834 -> 18624.75
$ws.Cells.Item(136,11).Value2 = 55874.25  # K136: 71699.50199999999 -> 55874.25
$ws.Cells.Item(136,13).Value2 = -53324.25  # M136: -69149.50199999999 -> -53324.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3,8).Value2 = 3328245  # H3: 3536229.8 -> 3328245
$ws.Cells.Item(3,9).Value2 = 5656372.5  # I3: 6284804 -> 5656372.5
$ws.Cells.Item(3,11).Value2 = 5656372.5  # K3: 6284804 -> 5656372.5
$ws.Cells.Item(3,13).Value2 = -5656258.5  # M3: -6284690 -> -5656258.5
$ws.Cells.Item(80,8).Value2 = 571.6923  # H80: 615.6923 -> 571.6923
$ws.Cells.Item(80,9).Value2 = 492  # I80: 596.5 -> 492
$ws.Cells.Item(80,10).Value2 = 607.1111  # J80: 619.1818 -> 607.1111
$ws.Cells.Item(80,11).Value2 = 492  # K80: 596.5 -> 492
$ws.Cells.Item(80,12).Value2 = 607.1111  # L80: 619.1818 -> 607.1111
$ws.Cells.Item(80,13).Value2 = 506  # M80: 401.5 -> 506
$ws.Cells.Item(80,14).Value2 = -2603.1111  # N80: -2615.1818 -> -2603.1111
$ws.Cells.Item(83,8).Value2 = 571.6923  # H83: 615.6923 -> 571.6923
$ws.Cells.Item(83,9).Value2 = 492  # I83: 596.5 -> 492
$ws.Cells.Item(83,10).Value2 = 607.1111  # J83: 619.1818 -> 607.1111
$ws.Cells.Item(83,11).Value2 = 2460  # K83: 2982.5 -> 2460
$ws.Cells.Item(83,12).Value2 = 3035.5555  # L83: 3095.909 -> 3035.5555
$ws.Cells.Item(83,13).Value2 = 2532  # M83: 2009.5 -> 2532
$ws.Cells.Item(83,14).Value2 = -13019.5555  # N83: -13079.909 -> -13019.5555
$ws.Cells.Item(94,8).Value2 = 2275983.5  # H94: 2328903.5 -> 2275983.5
$ws.Cells.Item(94,9).Value2 = 2501132  # I94: 2565252.5 -> 2501132
$ws.Cells.Item(94,11).Value2 = 2501132  # K94: 2565252.5 -> 2501132
$ws.Cells.Item(94,13).Value2 = -2500681  # M94: -2564801.5 -> -2500681
$ws.Cells.Item(132,8).Value2 = 0  # H132: 88999.5 -> 0
$ws.Cells.Item(132,10).Value2 = 0  # J132: 88999.5 -> 0
$ws.Cells.Item(132,12).ClearContents()  # L132: remove (was 88999.5)
$ws.Cells.Item(132,14).Value2 = 0  # N132: -99119.5 -> 0
$ws.Cells.Item(134,8).Value2 = 15278.454  # H134: 9794.085999999999 -> 15278.454
$ws.Cells.Item(134,9).Value2 = 12959.529  # I134: 7348.161 -> 12959.529
$ws.Cells.Item(134,10).Value2 = 23162.8  # J134: 28750 -> 23162.8
$ws.Cells.Item(134,11).Value2 = 38878.587  # K134: 22044.483 -> 38878.587
$ws.Cells.Item(134,12).Value2 = 69488.39999999999  # L134: 86250 -> 69488.39999999999
$ws.Cells.Item(134,13).Value2 = -36343.587  # M134: -19509.483 -> -36343.587
$ws.Cells.Item(134,14).Value2 = -74558.39999999999  # N134: -91320 -> -74558.39999999999
$ws.Cells.Item(135,8).Value2 = 129999  # H135: 122999.1 -> 129999
$ws.Cells.Item(135,10).Value2 = 129999  # J135: 122999.1 -> 129999
$ws.Cells.Item(135,12).Value2 = 129999  # L135: 122999.1 -> 129999
$ws.Cells.Item(135,14).Value2 = -140139  # N135: -133139.1 -> -140139

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value2 = 57054.473  # H31: 35339.242 -> 57054.473
$ws.Cells.Item(31,9).Value2 = 17055.285  # I31: 19606.334 -> 17055.285
$ws.Cells.Item(31,10).Value2 = 80387.336  # J31: 38835.445 -> 80387.336
$ws.Cells.Item(31,11).Value2 = 17055.285  # K31: 19606.334 -> 17055.285
$ws.Cells.Item(31,12).Value2 = 80387.336  # L31: 38835.445 -> 80387.336
$ws.Cells.Item(31,13).Value2 = -16760.285  # M31: -19311.334 -> -16760.285
$ws.Cells.Item(31,14).Value2 = -80977.336  # N31: -39425.445 -> -80977.336
$ws.Cells.Item(34,8).Value2 = 57054.473  # H34: 35339.242 -> 57054.473
$ws.Cells.Item(34,9).Value2 = 17055.285  # I34: 19606.334 -> 17055.285
$ws.Cells.Item(34,10).Value2 = 80387.336  # J34: 38835.445 -> 80387.336
$ws.Cells.Item(34,11).Value2 = 17055.285  # K34: 19606.334 -> 17055.285
$ws.Cells.Item(34,12).Value2 = 80387.336  # L34: 38835.445 -> 80387.336
$ws.Cells.Item(34,13).Value2 = -16853.285  # M34: -19404.334 -> -16853.285
$ws.Cells.Item(34,14).Value2 = -80791.336  # N34: -39239.445 -> -80791.336
$ws.Cells.Item(81,8).Value2 = 40000  # H81: 61500.5 -> 40000
$ws.Cells.Item(81,10).Value2 = 0  # J81: 83001 -> 0
$ws.Cells.Item(81,12).Value2 = 0  # L81: 83001 -> 0
$ws.Cells.Item(81,14).ClearContents()  # N81: remove (was -84997)
$ws.Cells.Item(84,8).Value2 = 40000  # H84: 61500.5 -> 40000
$ws.Cells.Item(84,10).Value2 = 0  # J84: 83001 -> 0
$ws.Cells.Item(84,12).Value2 = 0  # L84: 249003 -> 0
$ws.Cells.Item(84,14).ClearContents()  # N84: remove (was -258987)
$ws.Cells.Item(105,8).Value2 = 512.05  # H105: 512.55 -> 512.05
$ws.Cells.Item(105,9).Value2 = 407.83334  # I105: 408.3889 -> 407.83334
$ws.Cells.Item(105,11).Value2 = 407.83334  # K105: 408.3889 -> 407.83334
$ws.Cells.Item(105,13).Value2 = 1339.16666  # M105: 1338.6111 -> 1339.16666
$ws.Cells.Item(138,8).Value2 = 37991.332  # H138: 37991.5 -> 37991.332
$ws.Cells.Item(138,10).Value2 = 37991.332  # J138: 37991.5 -> 37991.332
$ws.Cells.Item(138,12).Value2 = 37991.332  # L138: 37991.5 -> 37991.332
$ws.Cells.Item(138,14).Value2 = -48271.332  # N138: -48271.5 -> -48271.332

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69,8).Value2 = 2800  # H69: 2000 -> 2800
$ws.Cells.Item(69,10).Value2 = 3000  # J69: 0 -> 3000
$ws.Cells.Item(69,12).Value2 = 9000  # L69: 0 -> 9000
$ws.Cells.Item(69,14).Value2 = -10622  # N69: None -> -10622
$ws.Cells.Item(72,8).Value2 = 2800  # H72: 2000 -> 2800
$ws.Cells.Item(72,10).Value2 = 3000  # J72: 0 -> 3000
$ws.Cells.Item(72,12).Value2 = 27000  # L72: 0 -> 27000
$ws.Cells.Item(72,14).Value2 = -35112  # N72: None -> -35112
$ws.Cells.Item(136,8).Value2 = 5462.8335  # H136: 5955.4 -> 5462.8335
$ws.Cells.Item(136,9).Value2 = 5555.4  # I136: 6194.25 -> 5555.4
$ws.Cells.Item(136,11).Value2 = 16666.2  # K136: 18582.75 -> 16666.2
$ws.Cells.Item(136,13).Value2 = -11566.2  # M136: -13482.75 -> -11566.2
$ws.Cells.Item(137,8).Value2 = 3790.3125  # H137: 3643.1333 -> 3790.3125
$ws.Cells.Item(137,10).Value2 = 7998.6665  # J137: 8999 -> 7998.6665
$ws.Cells.Item(137,12).Value2 = 23995.9995  # L137: 26997 -> 23995.9995
$ws.Cells.Item(137,14).Value2 = -34195.99950000001  # N137: -37197 -> -34195.99950000001
$ws.Cells.Item(141,8).Value2 = 2582.3845  # H141: 2524.5386 -> 2582.3845
$ws.Cells.Item(141,9).Value2 = 2568.4167  # I141: 2524.5386 -> 2568.4167
$ws.Cells.Item(141,10).Value2 = 2750  # J141: 0 -> 2750
$ws.Cells.Item(141,11).Value2 = 7705.250100000001  # K141: 7573.6158 -> 7705.250100000001
$ws.Cells.Item(141,12).Value2 = 8250  # L141: 0 -> 8250
$ws.Cells.Item(141,13).Value2 = -2525.250100000001  # M141: -2393.6158 -> -2525.250100000001
$ws.Cells.Item(141,14).Value2 = -18610  # N141: None -> -18610

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70,8).Value2 = 50003504  # H70: 66669668 -> 50003504
$ws.Cells.Item(70,10).Value2 = 5009  # J70: 0 -> 5009
$ws.Cells.Item(70,12).Value2 = 5009  # L70: 0 -> 5009
$ws.Cells.Item(70,14).Value2 = -5549  # N70: None -> -5549
$ws.Cells.Item(73,8).Value2 = 50003504  # H73: 66669668 -> 50003504
$ws.Cells.Item(73,10).Value2 = 5009  # J73: 0 -> 5009
$ws.Cells.Item(73,12).Value2 = 5009  # L73: 0 -> 5009
$ws.Cells.Item(73,14).Value2 = -6881  # N73: None -> -6881
$ws.Cells.Item(132,8).Value2 = 7266.6978  # H132: 7567.7075 -> 7266.6978
$ws.Cells.Item(132,9).Value2 = 5514.212  # I132: 5799.2583 -> 5514.212
$ws.Cells.Item(132,11).Value2 = 16542.636  # K132: 17397.7749 -> 16542.636
$ws.Cells.Item(132,13).Value2 = -14012.636  # M132: -14867.7749 -> -14012.636

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2,8).Value2 = 2879618.5  # H2: 2521541.2 -> 2879618.5
$ws.Cells.Item(2,9).Value2 = 10000000  # I2: 5007500 -> 10000000
$ws.Cells.Item(2,10).Value2 = 31466  # J2: 35582.75 -> 31466
$ws.Cells.Item(2,11).Value2 = 10000000  # K2: 5007500 -> 10000000
$ws.Cells.Item(2,12).Value2 = 31466  # L2: 35582.75 -> 31466
$ws.Cells.Item(2,13).Value2 = -9999888  # M2: -5007388 -> -9999888
$ws.Cells.Item(2,14).Value2 = -31690  # N2: -35806.75 -> -31690
$ws.Cells.Item(7,8).Value2 = 6059.091  # H7: 5403.4136 -> 6059.091
$ws.Cells.Item(7,9).Value2 = 4111.3076  # I7: 3842.3 -> 4111.3076
$ws.Cells.Item(7,11).Value2 = 4111.3076  # K7: 3842.3 -> 4111.3076
$ws.Cells.Item(7,13).Value2 = -3999.3076  # M7: -3730.3 -> -3999.3076
$ws.Cells.Item(126,8).Value2 = 6059.091  # H126: 5403.4136 -> 6059.091
$ws.Cells.Item(126,9).Value2 = 4111.3076  # I126: 3842.3 -> 4111.3076
$ws.Cells.Item(126,11).Value2 = 12333.9228  # K126: 11526.9 -> 12333.9228
$ws.Cells.Item(126,13).Value2 = -9863.9228  # M126: -9056.900000000001 -> -9863.9228

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4,8).Value2 = 145285.72  # H4: 146157.14 -> 145285.72
$ws.Cells.Item(4,9).Value2 = 460499.5  # I4: 920000 -> 460499.5
$ws.Cells.Item(4,10).Value2 = 19200.2  # J4: 17183.334 -> 19200.2
$ws.Cells.Item(4,11).Value2 = 460499.5  # K4: 920000 -> 460499.5
$ws.Cells.Item(4,12).Value2 = 19200.2  # L4: 17183.334 -> 19200.2
$ws.Cells.Item(4,13).Value2 = -460386.5  # M4: -919887 -> -460386.5
$ws.Cells.Item(4,14).Value2 = -19426.2  # N4: -17409.334 -> -19426.2
$ws.Cells.Item(132,8).Value2 = 14657606  # H132: 15323848 -> 14657606
$ws.Cells.Item(132,9).Value2 = 18188426  # I132: 19237742 -> 18188426
$ws.Cells.Item(132,11).Value2 = 54565278  # K132: 57713226 -> 54565278
$ws.Cells.Item(132,13).Value2 = -54562748  # M132: -57710696 -> -54562748
